$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 for 18.9.2025 working hours
$ws.Range("A10").Value = "18.9.2025"

$ws.Range("B10").Value = 0.5
$ws.Range("C10").Value = 0.60416666666666663
$ws.Range("D10").Value = 0.88541666666666663
$ws.Range("E10").Value = 0.90972222222222221

# Match the time number format used by the other rows (style index 1 -> numFmtId 18, h:mm AM/PM)
$ws.Range("B10:E10").NumberFormat = $ws.Range("B9:E9").NumberFormat
